$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J23").Select() | Out-Null

$ws.Range("A98").Value = "fuzzy match ranking"
$ws.Range("A83").Value = "konvertuj utc vremena u epoh vreme"
$ws.Range("A74").Value = "toplotna mapa 3D koordinata"
$ws.Range("A73").Value = "ekstrakcija podatka iz html sadrzaja"
$ws.Range("A75").Value = "uzeti sve roditelje od xml cvora"
$ws.Range("A77").Value = "podvući tekst u nazivu vidzeta"
$ws.Range("A88").Value = "pamcenje na disku - stalno skladiste"
$ws.Range("A90").Value = "kako citati sadrzaj .gz kompresovanog fajla?"
$ws.Range("A100").Value = "kako čitati .csv fajl na efikasan nacin?"
$ws.Range("A70").Value = "mnozenje matrica"
